# TC31_Verify_login.xlsx -- "Changes for New UI Prod"
#
# This script reproduces the data/structure edits captured by the commit:
#   - TC31_Verify_login sheet: swap some CLICK/VERIFY steps for MOUSEOVER,
#     rename the post-login verification target from "Welcomeelement" to
#     "MyaccountSection", retarget the final CLICK to "Logout" and drop the
#     now-redundant trailing row, tweak a couple of row heights, and move
#     the active selection.
#   - Testdata sheet: add two new data rows (EleType1 / EleType2, both
#     backed by "JSElement").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TC31_Verify_login")
$ws2 = $wb.Worksheets.Item("Testdata")

# ---------------------------------------------------------------------
# Sheet 1: TC31_Verify_login
# ---------------------------------------------------------------------

# Row 3: LoginOption click becomes a mouseover
$ws1.Cells.Item(3, 2).Value = "MOUSEOVER"

# Row 9: the element verified after login is now "MyaccountSection"
# (previously "Welcomeelement")
$ws1.Cells.Item(9, 3).Value = "MyaccountSection"
$ws1.Cells.Item(9, 5).Value = "MyaccountSection"

# Row 10: verifying MyaccountSection becomes a mouseover step, and the
# data-descriptor column is no longer populated
$ws1.Cells.Item(10, 2).Value = "MOUSEOVER"
$ws1.Cells.Item(10, 5).ClearContents()

# Row 11: the final click now targets "Logout" directly
$ws1.Cells.Item(11, 3).Value = "Logout"

# Row 12 (old standalone "CLICK Logout" row) is no longer needed now that
# row 11 covers it - remove the whole row
$ws1.Rows.Item(12).Delete()

# A couple of rows picked up explicit custom heights
$ws1.Rows.Item(6).RowHeight = 15.75
$ws1.Rows.Item(8).RowHeight = 14.25

# Move the active selection
$ws1.Range("A9:XFD9").Select()

# ---------------------------------------------------------------------
# Sheet 2: Testdata
# ---------------------------------------------------------------------

# New data rows backing the EleType1 / EleType2 descriptors used above
$ws2.Cells.Item(8, 1).Value = "EleType1"
$ws2.Cells.Item(8, 2).Value = "JSElement"
$ws2.Cells.Item(9, 1).Value = "EleType2"
$ws2.Cells.Item(9, 2).Value = "JSElement"

# Give the two new rows the same thin-border look used by the rest of the table
$ws2.Range("A8:B9").Borders.LineStyle = 1
